# Applies the "adding descriptions for fishery shenanigans" edit:
#  1. Inserts two new narrative paragraphs (plus a bottom-bordered blank
#     paragraph and two plain blank paragraphs) right before the
#     "Thinking about the Pinniped Components" heading paragraph.
#  2. Stamps a <w:lastRenderedPageBreak/> onto the run that starts the
#     "I think baseline take..." paragraph.
#
# Note: this COM shim's Paragraph.Previous / Paragraph.Next navigation
# returns objects whose Range has no resolvable Start/End, so we use
# Paragraphs.Item(index) (1-based) throughout instead.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: insert the new fishery-data paragraphs.
# ---------------------------------------------------------------------

$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    if ($p.Range.Text -match "^Thinking about the Pinniped Components") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate 'Thinking about the Pinniped Components' paragraph"
}

$anchorIndex = $targetIndex - 1
$anchor = $d.Paragraphs.Item($anchorIndex)

# Create one fresh paragraph right after the anchor; InsertXML below will
# expand it into the full run of new paragraphs while leaving everything
# else in the document untouched.
$anchor.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 1)

# Note: InsertXML's trailing paragraph mark is absorbed into the target
# (pre-existing) paragraph mark rather than creating a new paragraph, so
# the fragment below carries one extra trailing <w:p/> to land on the
# intended five new paragraphs (two narrative + bordered blank + two
# plain blanks).
$fisheryXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Fishery data is in sheet &#x201C;Chinook and Chum catches for Liz.xlsx&#x201D; including boat counts and landings and effort. But we went back and decided on a more deliberate time frame for fishery data and that&#x2019;s captured in &#x201C;</w:t></w:r><w:r><w:t>Adjusted_Nisqually_Chinook_and_Chum_from_Craig_July2024 and August.xlsx"</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Chum fishery is somewhere between 1-25 boats in any year, Chinook 10-20 estimated. </w:t></w:r><w:r><w:t>Average boats during chum 2014 &#x2013; 2019 = 12. Average boats during Chinook/Coho is very rough estimate from Craig, about 14, 2017-2023.</w:t></w:r></w:p><w:p><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr></w:pPr></w:p><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newPara.Range.InsertXML($fisheryXml) | Out-Null

# ---------------------------------------------------------------------
# Part 2: add <w:lastRenderedPageBreak/> to the "I think baseline take"
# paragraph's run, preserving its existing paragraph identity.
# ---------------------------------------------------------------------

$baselineIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    if ($p.Range.Text -match "^I think baseline take might look like") {
        $baselineIndex = $i
        break
    }
}

if ($baselineIndex -eq -1) {
    throw "Could not locate the 'I think baseline take...' paragraph"
}

$baseline = $d.Paragraphs.Item($baselineIndex)

# Replace only the paragraph's text content (excluding its trailing
# paragraph mark) so the <w:p> element keeps its original w14:paraId /
# rsid* attributes intact - only the run itself changes.
$bStart = $baseline.Range.Start
$bEnd = $baseline.Range.End - 1
$bodyRange = $d.Range($bStart, $bEnd)

$pageBreakXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>I think baseline take might look like a low level of take on all fishery opener days. Ask about species specific rates/preferences?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$bodyRange.InsertXML($pageBreakXml) | Out-Null

Write-Output "done"
